$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.5
$ws.Range("I3").Value = 2.6
$ws.Range("U3").Value = 13
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 23
$ws.Range("Z3").Value = 13
$ws.Range("AE3").Value = 11
$ws.Range("AF3").Value = 15
$ws.Range("AI3").Value = 21

# Row 6
$ws.Range("N6").Value = 1.98
$ws.Range("O6").Value = 1.88

# Row 7
$ws.Range("J7").Value = 1.05
$ws.Range("K7").Value = 11
$ws.Range("L7").Value = 1.29
$ws.Range("M7").Value = 3.5
$ws.Range("N7").Value = 1.9
$ws.Range("O7").Value = 1.95

# Row 8
$ws.Range("N8").Value = 1.85
$ws.Range("O8").Value = 2

# Row 12
$ws.Range("G12").Value = 3.5
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 2.15
$ws.Range("J12").Value = 1.1
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = 1.44
$ws.Range("M12").Value = 2.63
$ws.Range("N12").Value = 2.5
$ws.Range("P12").Value = 1.53
$ws.Range("Q12").Value = 2.38
$ws.Range("R12").Value = 2.1
$ws.Range("S12").Value = 1.67
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 17
$ws.Range("V12").Value = 13
$ws.Range("W12").Value = 41
$ws.Range("X12").Value = 34
$ws.Range("Y12").Value = 41
$ws.Range("Z12").Value = 7
$ws.Range("AA12").Value = 6
$ws.Range("AC12").Value = 67
$ws.Range("AF12").Value = 9
$ws.Range("AG12").Value = 10
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 21
$ws.Range("AJ12").Value = 41

# Row 13
$ws.Range("T13").Value = 6.3

# Row 15
$ws.Range("G15").Value = 2.7
$ws.Range("I15").Value = 2.52
$ws.Range("L15").Value = 1.37
$ws.Range("M15").Value = 2.62
$ws.Range("N15").Value = 2.07
$ws.Range("Q15").Value = 2.45
$ws.Range("T15").Value = 7.8
$ws.Range("U15").Value = 13
$ws.Range("V15").Value = 10.25
$ws.Range("W15").Value = 32
$ws.Range("X15").Value = 25
$ws.Range("Y15").Value = 37
$ws.Range("AA15").Value = 6
$ws.Range("AE15").Value = 7.2
$ws.Range("AF15").Value = 11.75
$ws.Range("AG15").Value = 9.75
$ws.Range("AH15").Value = 27
$ws.Range("AI15").Value = 23
$ws.Range("AJ15").Value = 37

# Row 23
$ws.Range("G23").Value = 2
$ws.Range("I23").Value = 4.1
$ws.Range("J23").Value = 1.08
$ws.Range("K23").Value = 8
$ws.Range("N23").Value = 2.1
$ws.Range("O23").Value = 1.7
$ws.Range("R23").Value = 1.91
$ws.Range("S23").Value = 1.91
$ws.Range("U23").Value = 9
$ws.Range("Z23").Value = 8
$ws.Range("AD23").Value = 301
$ws.Range("AF23").Value = 19

# Row 24
$ws.Range("G24").Value = 1.85
$ws.Range("I24").Value = 5
$ws.Range("J24").Value = 1.07
$ws.Range("K24").Value = 8.5
$ws.Range("L24").Value = 1.3
$ws.Range("M24").Value = 3.4
$ws.Range("P24").Value = 1.4
$ws.Range("Q24").Value = 2.75
$ws.Range("X24").Value = 15
$ws.Range("Z24").Value = 8.5
$ws.Range("AF24").Value = 23

# Row 26
$ws.Range("J26").Value = 1.05
$ws.Range("K26").Value = 11

# Row 29
$ws.Range("R29").Value = 2.16
$ws.Range("S29").Value = 1.62

# Row 30
$ws.Range("G30").Value = 1.7
$ws.Range("H30").Value = 3.45
$ws.Range("I30").Value = 4.25
$ws.Range("L30").Value = 1.29
$ws.Range("M30").Value = 3.25
$ws.Range("N30").Value = 1.85
$ws.Range("O30").Value = 1.75
$ws.Range("Q30").Value = 2.47
$ws.Range("T30").Value = 5.8
$ws.Range("V30").Value = 6.9
$ws.Range("W30").Value = 11
$ws.Range("X30").Value = 11.25
$ws.Range("Y30").Value = 21
$ws.Range("Z30").Value = 9.5
$ws.Range("AA30").Value = 5.9
$ws.Range("AC30").Value = 55
$ws.Range("AF30").Value = 19
$ws.Range("AG30").Value = 12
$ws.Range("AH30").Value = 55
$ws.Range("AI30").Value = 32
$ws.Range("AJ30").Value = 37
